# Add two new columns "I0" (column I) and "IF" (column J) to the sheet.
# Column I0 is always 1; column IF mirrors the existing IP column (H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 19

# Copy the header formatting from H1 (bold, centered, bordered) onto the
# new header cells I1:J1 before setting their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
